$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows to match re-pulled data
$ws.Range("F2").Value = -2
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = 0
$ws.Range("F9").Value = 6
$ws.Range("F10").Value = 0
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = -4
$ws.Range("F15").Value = 4
$ws.Range("F16").Value = -5
$ws.Range("F22").Value = -7
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 3
